$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, reporting readings for the cities queried by the user
$data = @(
    @("16/02/2025 21:42:34", "São Paulo", "27", "58%", "Alerta Amarelo, Média Umidade no ar"),
    @("16/02/2025 21:46:05", "São Paulo", "27", "58%", "Alerta Amarelo, Média Umidade no ar"),
    @("16/02/2025 21:46:57", "Taboão da Serra", "28", "60%", "Alerta Amarelo, Média Umidade no ar")
)

$startRow = 7
$endRow = $startRow + $data.Length - 1
$rng = $ws.Range("A$($startRow):E$($endRow)")

# Force text storage so values like "27" and "58%" are not reinterpreted as
# numbers/percentages by Excel, matching the plain string cells used
# elsewhere in the sheet.
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# Restore the default "Normal" style so the new cells don't carry an
# explicit text number format style, keeping them consistent with the
# rest of the sheet (the values remain text since they were entered as
# strings).
$rng.Style = "Normal"
